$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = 160

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value = $newVal
        }
    }
}
